$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 7: new data row inserted into the first table (mirrors rows 2-6)
# ------------------------------------------------------------------
$ws.Range("G7").Value = 2912000
$ws.Range("H7").Value = 14
$ws.Range("I7").Value = 14

# Copy the "Berechnung" cell format (with the thin box border) from row 6
# down onto the new J7 / M7 / O7 cells before writing their formulas.
$ws.Range("J6").Copy() | Out-Null
$ws.Range("J7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("M6").Copy() | Out-Null
$ws.Range("M7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("O6").Copy() | Out-Null
$ws.Range("O7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("J7").Formula = "=G7/H7/I7"
$ws.Range("M7").Formula = "=364/H7"
$ws.Range("O7").Formula = "=J7/M7"

# ------------------------------------------------------------------
# Row 15: new data row extending the second table (A13:D14 -> A15:D15)
# ------------------------------------------------------------------
$ws.Range("A15").Value = 20000
$ws.Range("B15").Value = 14
$ws.Range("C15").Value = 14

# Copy the format from D14 (same "Berechnung" style) down onto D15.
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("D15").Formula = "=A15*B15*C15"

# ------------------------------------------------------------------
# Selection follows the newly entered cell, like in the source edit.
# ------------------------------------------------------------------
$ws.Range("D15").Select() | Out-Null
